$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 holds the sample student record: update name + phone number,
# keeping the phone number stored as text (leading apostrophe prevents
# Excel from re-interpreting the digit string as a numeric value).
$ws.Range("A3").Value = "Adek"
$ws.Range("B3").Value = "'62895396334564"
